$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-03-15 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-16 Sunday", 2) | Out-Null
$d.Content.Find.Execute("88×42=3696", $true, $false, $false, $false, $false, $true, 1, $false, "99×88=8712", 2) | Out-Null
$d.Content.Find.Execute("23×68=1564", $true, $false, $false, $false, $false, $true, 1, $false, "15×45=675", 2) | Out-Null
$d.Content.Find.Execute("22×46=1012", $true, $false, $false, $false, $false, $true, 1, $false, "89×27=2403", 2) | Out-Null
$d.Content.Find.Execute("21×72=1512", $true, $false, $false, $false, $false, $true, 1, $false, "13×17=221", 2) | Out-Null
$d.Content.Find.Execute("75×28=2100", $true, $false, $false, $false, $false, $true, 1, $false, "91×52=4732", 2) | Out-Null
$d.Content.Find.Execute("21×42=882", $true, $false, $false, $false, $false, $true, 1, $false, "66×74=4884", 2) | Out-Null
$d.Content.Find.Execute("75×77=5775", $true, $false, $false, $false, $false, $true, 1, $false, "99×83=8217", 2) | Out-Null
$d.Content.Find.Execute("15×89=1335", $true, $false, $false, $false, $false, $true, 1, $false, "96×88=8448", 2) | Out-Null
$d.Content.Find.Execute("91×96=8736", $true, $false, $false, $false, $false, $true, 1, $false, "18×68=1224", 2) | Out-Null
$d.Content.Find.Execute("57×16=912", $true, $false, $false, $false, $false, $true, 1, $false, "34×47=1598", 2) | Out-Null
$d.Content.Find.Execute("84×14=1176", $true, $false, $false, $false, $false, $true, 1, $false, "93×64=5952", 2) | Out-Null
$d.Content.Find.Execute("62×31=1922", $true, $false, $false, $false, $false, $true, 1, $false, "31×18=558", 2) | Out-Null
$d.Content.Find.Execute("95×17=1615", $true, $false, $false, $false, $false, $true, 1, $false, "84×26=2184", 2) | Out-Null
$d.Content.Find.Execute("48×33=1584", $true, $false, $false, $false, $false, $true, 1, $false, "78×34=2652", 2) | Out-Null
$d.Content.Find.Execute("13×68=884", $true, $false, $false, $false, $false, $true, 1, $false, "83×40=3320", 2) | Out-Null
$d.Content.Find.Execute("38×62=2356", $true, $false, $false, $false, $false, $true, 1, $false, "18×51=918", 2) | Out-Null
$d.Content.Find.Execute("99×60=5940", $true, $false, $false, $false, $false, $true, 1, $false, "88×67=5896", 2) | Out-Null
$d.Content.Find.Execute("63×31=1953", $true, $false, $false, $false, $false, $true, 1, $false, "41×26=1066", 2) | Out-Null
$d.Content.Find.Execute("50×59=2950", $true, $false, $false, $false, $false, $true, 1, $false, "29×45=1305", 2) | Out-Null
$d.Content.Find.Execute("97×62=6014", $true, $false, $false, $false, $false, $true, 1, $false, "28×45=1260", 2) | Out-Null
$d.Content.Find.Execute("34×59=2006", $true, $false, $false, $false, $false, $true, 1, $false, "35×45=1575", 2) | Out-Null
$d.Content.Find.Execute("19×21=399", $true, $false, $false, $false, $false, $true, 1, $false, "73×19=1387", 2) | Out-Null
$d.Content.Find.Execute("68×15=1020", $true, $false, $false, $false, $false, $true, 1, $false, "29×13=377", 2) | Out-Null
$d.Content.Find.Execute("53×57=3021", $true, $false, $false, $false, $false, $true, 1, $false, "40×78=3120", 2) | Out-Null
$d.Content.Find.Execute("53×58=3074", $true, $false, $false, $false, $false, $true, 1, $false, "37×66=2442", 2) | Out-Null
